$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $searchStart = 0
    $docEnd = $d.Content.End
    $iterations = 0
    while ($searchStart -lt $docEnd -and $iterations -lt 1000) {
        $iterations = $iterations + 1
        $rng = $d.Range($searchStart, $docEnd)
        $rng.Find.ClearFormatting()
        $rng.Find.Text = $old
        $rng.Find.Forward = $true
        $rng.Find.Wrap = 0
        $rng.Find.MatchCase = $true
        $rng.Find.MatchWholeWord = $false
        $rng.Find.MatchWildcards = $false
        $rng.Find.Execute() | Out-Null
        if (-not $rng.Find.Found) {
            break
        }
        $rng.Text = $new
        $searchStart = $rng.End
        $docEnd = $d.Content.End
    }
}

Replace-All "1122334" "1231231"
Replace-All "07.03.2024" "06.06.2024"
Replace-All "Акционерное общество «Арктические морские инженерно-геологические экспедиции»" "Атлантическая база флота - филиал ФГБУ науки Института океанологии им. П.П. Ширшова Российской академии наук"
Replace-All "Капитан Бахтин Ю. Г." "заместитель директора безопасности Волков А. В."
Replace-All "Кодекса торгового мореплавания (КТМ РФ)" "Доверенности № 1 от 01.06.2024"
Replace-All '"СИНЕГОРСК" ' '"АКАДЕМИК ИОФФЕ" '
Replace-All "021026" "870072"
Replace-All "Первоначальное освидетельствование" "Ежегодное освидетельствование"
Replace-All "Свидетельство ф. 8.5.3 № 24.42.02.00123.121 от 04.05.2024" "Акт ф. 6.1.03 № 123213213123 от --"
Replace-All "100 000,00 p. (сто тысяч рублей 00 копеек)" "123 123,00 p. (сто двадцать три тысячи сто двадцать три рубля 00 копеек)"
Replace-All "Ю. Г. Бахтин" "А. В. Волков"
